$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The FilesTab Neo4j query (cell B4) was corrected: the `File Type` and
# `Breed` columns were dropped from the RETURN clause (and a couple of
# lines re-indented) as part of "corrected ICDC Breed 1-14 scripts".
$newFilesTabQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Labrador Retriever']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
         coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesTabQuery

# Match the refreshed selection/view state left behind after the edit.
$ws.Range("B4").Select()

$wb.Save()
